$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.498.11"
$ws.Range("E2").Value = "  +2.78%  "
$ws.Range("D3").Value = "1.603.70"
$ws.Range("E3").Value = "  +2.36%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("E6").Value = "  +6.73%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "26.79"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.52"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.20%  "
$ws.Range("E10").Value = "  +2.49%  "
$ws.Range("E11").Value = "  +2.69%  "
$ws.Range("E12").Value = "  +1.63%  "
$ws.Range("D13").Value = "1.833.73"
$ws.Range("E13").Value = "  +2.44%  "
$ws.Range("D14").Value = "1.639.08"
$ws.Range("E14").Value = "  +4.53%  "
$ws.Range("D15").Value = "29.485.50"
$ws.Range("E15").Value = "  +2.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.536"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.43%  "
$ws.Range("E17").Value = "  +1.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.42%  "
$ws.Range("D21").Value = "0.0₃0692"
$ws.Range("E21").Value = "  +1.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  +1.70%  "
$ws.Range("E24").Value = "  +1.84%  "
$ws.Range("E25").Value = "  +0.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.67%  "
$ws.Range("E27").Value = "  +4.95%  "
$ws.Range("E28").Value = "  +3.14%  "
$ws.Range("E29").Value = "  +2.29%  "
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("E31").Value = "  +2.58%  "
$ws.Range("E32").Value = "  +0.98%  "
$ws.Range("E33").Value = "  +1.61%  "
$ws.Range("E34").Value = "  +3.45%  "
$ws.Range("D35").Value = "1.409.44"
$ws.Range("E35").Value = "  +0.96%  "
$ws.Range("E36").Value = "  +1.04%  "
$ws.Range("E37").Value = "  +3.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.82"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.85%  "
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0165"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.42%  "
$ws.Range("E41").Value = "  +3.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.98"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0486"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "52.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +22.54%  "
$ws.Range("E45").Value = "  +3.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.998"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "65.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.68%  "
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("D49").Value = "1.742.93"
$ws.Range("E49").Value = "  +2.07%  "
$ws.Range("E50").Value = "  -1.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "86.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.54%  "
